$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.231.79'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.861.36'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7160'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '240.77'
$ws.Range('E6').Value = '  +0.57%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07759'
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3084'
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.16'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08262'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.240'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.850.04'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.7177'
$ws.Range('E14').Value = '  -0.81%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.31'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.208.82'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.871'
$ws.Range('E17').Value = '  +0.82%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '244.55'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007815'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.17'
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.109.19'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.956'
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1593'
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '162.82'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.944'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.28'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.496'
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('E30').Value = '  -3.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.407'
$ws.Range('E31').Value = '  +1.63%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.171'
$ws.Range('E32').Value = '  +2.39%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05204'
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.912'
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.174'
$ws.Range('E35').Value = '  -2.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7287'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01856'
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.689'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.154.37'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9050'
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.101'
$ws.Range('E42').Value = '  +1.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '72.50'
$ws.Range('E43').Value = '  +1.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '101.78'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.004.48'
$ws.Range('E46').Value = '  -1.71%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5237'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.770'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.340'
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.872'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.068'
$ws.Range('E51').Value = '  +0.32%  '
